$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be stored as text so numeric-looking price strings
# (e.g. "234.91") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.252.32"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.870.89"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "234.91"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "0.4698"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "0.2852"
$ws.Range("E8").Value = "  -1.50%  "

$ws.Range("D9").Value = "41.65"
$ws.Range("E9").Value = "  -2.82%  "

$ws.Range("D10").Value = "0.06551"
$ws.Range("E10").Value = "  +0.45%  "

$ws.Range("D11").Value = "21.24"
$ws.Range("E11").Value = "  -1.79%  "

$ws.Range("D12").Value = "0.07825"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").Value = "96.74"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").Value = "1.851.77"
$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("D15").Value = "0.6921"
$ws.Range("E15").Value = "  +2.75%  "

$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "268.43"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "30.248.46"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").Value = "13.78"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").Value = "0.000007691"
$ws.Range("E20").Value = "  +2.01%  "

$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").Value = "2.140.73"
$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "5.249"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").Value = "6.172"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").Value = "9.492"
$ws.Range("E26").Value = "  +3.69%  "

$ws.Range("D27").Value = "166.05"
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("D28").Value = "18.81"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "1.935"
$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("D30").Value = "1.371"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("D31").Value = "0.09939"
$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("D32").Value = "4.357"
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").Value = "1.455"
$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").Value = "4.058"
$ws.Range("E34").Value = "  +1.66%  "

$ws.Range("D35").Value = "0.04748"
$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("D36").Value = "1.131"
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").Value = "0.7015"
$ws.Range("E37").Value = "  +0.89%  "

$ws.Range("D38").Value = "2.719"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("D40").Value = "2.779"
$ws.Range("E40").Value = "  +6.83%  "

$ws.Range("D41").Value = "6.282"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D42").Value = "73.00"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "1.939"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").Value = "0.4165"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D48").Value = "982.84"
$ws.Range("E48").Value = "  +4.49%  "

$ws.Range("D49").Value = "7.105"
$ws.Range("E49").Value = "  +2.20%  "

$ws.Range("D50").Value = "9.178"
$ws.Range("E50").Value = "  +0.74%  "

$ws.Range("D51").Value = "34.53"
$ws.Range("E51").Value = "  +2.49%  "

# Row 46 and 47: coin order swapped (Quant now ranked 44, TrustWalletToken 45)
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "103.21"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "0.8341"
$ws.Range("E47").Value = "  -0.28%  "
